# Weekly driver report update for 2025-04-21
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

# Row 3: Intel(R) Wireless-AC 9260 160MHz - 23.40.0.4
$ws.Range("C3").Value = 367
$ws.Range("D3").Value = 97.40000000000001

# Rows 4 and 5 swap places (with refreshed counts):
# Row 4 becomes: Intel(R) Wi-Fi 6 AX200 160MHz - 23.120.0.3
$ws.Range("A4").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 23.120.0.3"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 98
$ws.Range("D4").Value = 98.3

# Row 5 becomes: Intel(R) Dual Band Wireless-AC 8265 - 20.70.25.2
$ws.Range("A5").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.25.2"
$ws.Range("B5").Value = 8
$ws.Range("C5").Value = 1143
$ws.Range("D5").Value = 98.5

# Row 6: Intel(R) Wi-Fi 6 AX200 160MHz - 22.250.0.4
$ws.Range("C6").Value = 1011

# Row 7: Totals
$ws.Range("C7").Value = 2619

# Row 24: Intel(R) Wi-Fi 6 AX200 160MHz - 22.230.0.8
$ws.Range("B24").Value = 331283

# Row 42: Intel(R) Dual Band Wireless-AC 8265 - 20.70.12.5
$ws.Range("B42").Value = 144782
